$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "37.430.54"
$ws.Range("E2").Value = "  +5.51%  "

# Row 3
$ws.Range("D3").Value = "2.037.54"
$ws.Range("E3").Value = "  +3.05%  "

# Row 4
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").Value = "'253.08"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +5.91%  "

# Row 6
$ws.Range("D6").Value = "'0.645"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.96%  "

# Row 7
$ws.Range("D7").Value = "'62.95"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +14.04%  "

# Row 8
$ws.Range("E8").Value = "  +0.17%  "

# Row 9
$ws.Range("D9").Value = "'0.374"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +6.61%  "

# Row 10
$ws.Range("D10").Value = "'58.93"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.17%  "

# Row 11
$ws.Range("D11").Value = "'0.0753"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +4.74%  "

# Row 12
$ws.Range("E12").Value = "  +0.42%  "

# Row 13
$ws.Range("D13").Value = "'0.911"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +3.92%  "

# Row 14
$ws.Range("D14").Value = "'15.03"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +7.66%  "

# Row 15
$ws.Range("D15").Value = "2.340.13"
$ws.Range("E15").Value = "  +3.31%  "

# Row 16
$ws.Range("D16").Value = "'5.56"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +7.57%  "

# Row 17
$ws.Range("D17").Value = "'20.46"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +20.46%  "

# Row 18
$ws.Range("D18").Value = "2.051.73"
$ws.Range("E18").Value = "  +3.74%  "

# Row 19
$ws.Range("D19").Value = "37.382.95"
$ws.Range("E19").Value = "  +5.81%  "

# Row 20
$ws.Range("D20").Value = "'73.17"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +5.20%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0870"
$ws.Range("E21").Value = "  +5.34%  "

# Row 22
$ws.Range("D22").Value = "'5.32"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +7.72%  "

# Row 23
$ws.Range("D23").Value = "'236.14"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.55%  "

# Row 24
$ws.Range("D24").Value = "'2.76"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +24.93%  "

# Row 25
$ws.Range("E25").Value = "  -0.16%  "

# Row 26
$ws.Range("D26").Value = "'2.32"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.80%  "

# Row 27
$ws.Range("D27").Value = "'9.50"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +5.73%  "

# Row 28
$ws.Range("D28").Value = "'165.49"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.19%  "

# Row 29
$ws.Range("D29").Value = "'19.78"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.89%  "

# Row 30
$ws.Range("D30").Value = "'0.121"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.44%  "

# Row 31
$ws.Range("D31").Value = "'0.113"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +28.43%  "

# Row 32
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'5.17"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +10.00%  "

# Row 33
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'1.21"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +6.98%  "

# Row 34
$ws.Range("D34").Value = "'4.69"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +11.78%  "

# Row 35
$ws.Range("D35").Value = "'0.0613"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +5.78%  "

# Row 36
$ws.Range("D36").Value = "'2.43"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +12.81%  "

# Row 37
$ws.Range("E37").Value = "  -0.08%  "

# Row 38
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").Value = "'1.81"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.31%  "

# Row 39
$ws.Range("B39").Value = "THORChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D39").Value = "'5.99"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +25.01%  "

# Row 40
$ws.Range("D40").Value = "'0.104"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +20.04%  "

# Row 41
$ws.Range("D41").Value = "'1.23"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +6.18%  "

# Row 42
$ws.Range("E42").Value = "  +4.19%  "

# Row 43
$ws.Range("E43").Value = "  +6.17%  "

# Row 44
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'2.71"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +22.87%  "

# Row 45
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").Value = "'1.14"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +6.41%  "

# Row 46
$ws.Range("D46").Value = "'8.07"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +10.90%  "

# Row 47
$ws.Range("D47").Value = "'16.85"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +10.68%  "

# Row 48
$ws.Range("D48").Value = "'94.97"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +7.04%  "

# Row 49
$ws.Range("D49").Value = "1.422.70"
$ws.Range("E49").Value = "  +4.39%  "

# Row 50
$ws.Range("D50").Value = "'2.94"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.11%  "

# Row 51
$ws.Range("D51").Value = "'47.34"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +6.30%  "
